# Generate Report for Handoff
# Updates the localization-status report:
#  - Overview sheet: refresh "Latest HO Xliff Generate Date" for the rows that
#    were just (re)handed off
#  - de-de sheet: same handoff timestamp refresh for its "Latest Handoff
#    Datetime" column
#  - zh-cn sheet: its own (later) handoff timestamp refresh
#  - zh-cn / de-de sheets: set the "Priority" column to "ht" for those rows

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 13)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-28 18:22:37"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-28 18:22:32"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-28 18:22:37"
}
